$wb = $excel.ActiveWorkbook

# --- Sheet: ACCU ---
$ws = $wb.Worksheets.Item("ACCU")
$ws.Range("B2").Value = -986487.1946236491
$ws.Range("D2").Value = 544804.4168332688
$ws.Range("E2").Value = -441682.7777903803
$ws.Range("B3").Value = -917981.1394414512
$ws.Range("D3").Value = 479962.8082025978
$ws.Range("E3").Value = -438018.3312388534
$ws.Range("B4").Value = -849475.0842592533
$ws.Range("D4").Value = 417337.3845484682
$ws.Range("E4").Value = -432137.6997107852
$ws.Range("B5").Value = -780969.0290770555
$ws.Range("D5").Value = 357645.3585419968
$ws.Range("E5").Value = -423323.6705350587
$ws.Range("B6").Value = -712462.9738948577
$ws.Range("D6").Value = 301589.6717097988
$ws.Range("E6").Value = -410873.3021850589
$ws.Range("B7").Value = -643956.9187126597
$ws.Range("D7").Value = 249788.0716214628
$ws.Range("E7").Value = -394168.8470911969
$ws.Range("B8").Value = -575450.8635304619
$ws.Range("D8").Value = 202714.1631296072
$ws.Range("E8").Value = -372736.7004008548
$ws.Range("B9").Value = -506944.8083482641
$ws.Range("D9").Value = 160659.70250045
$ws.Range("E9").Value = -346285.1058478141
$ws.Range("B10").Value = -438438.7531660662
$ws.Range("D10").Value = 123721.6371589445
$ws.Range("E10").Value = -314717.1160071217
$ws.Range("B11").Value = -369932.6979838684
$ws.Range("D11").Value = 91812.04858423567
$ws.Range("E11").Value = -278120.6493996327
$ws.Range("B12").Value = -301426.6428016705
$ws.Range("D12").Value = 64685.4942794268
$ws.Range("E12").Value = -236741.1485222437
$ws.Range("B13").Value = -232920.5876194726
$ws.Range("D13").Value = 41976.65671461845
$ws.Range("E13").Value = -190943.9309048541
$ws.Range("B14").Value = -164414.5324372748
$ws.Range("D14").Value = 23241.40862297486
$ws.Range("E14").Value = -141173.1238142999
$ws.Range("B15").Value = -95908.47725507691
$ws.Range("D15").Value = 7995.770728413785
$ws.Range("E15").Value = -87912.70652666313
$ws.Range("B16").Value = -27402.42207287905
$ws.Range("D16").Value = -4250.89667133622
$ws.Range("E16").Value = -31653.31874421527
$ws.Range("B17").Value = 41103.63310931881
$ws.Range("D17").Value = -13970.28805257787
$ws.Range("E17").Value = 27133.34505674095
$ws.Range("B18").Value = 109609.6882915167
$ws.Range("D18").Value = -21597.45791055613
$ws.Range("E18").Value = 88012.23038096054
$ws.Range("B19").Value = 178115.7434737145
$ws.Range("D19").Value = -27520.33610708088
$ws.Range("E19").Value = 150595.4073666336
$ws.Range("B20").Value = 246621.7986559124
$ws.Range("D20").Value = -32075.18343470259
$ws.Range("E20").Value = 214546.6152212098
$ws.Range("B21").Value = 315127.8538381102
$ws.Range("D21").Value = -35546.5167697529
$ws.Range("E21").Value = 279581.3370683573
$ws.Range("B22").Value = 383633.9090203081
$ws.Range("D22").Value = -38170.08659593859
$ws.Range("E22").Value = 345463.8224243695
$ws.Range("B23").Value = 452139.9642025059
$ws.Range("D23").Value = -40137.68958333765
$ws.Range("E23").Value = 412002.2746191683
$ws.Range("B24").Value = 520646.0193847038
$ws.Range("D24").Value = -41602.86029687145
$ws.Range("E24").Value = 479043.1590878324
$ws.Range("B25").Value = 589152.0745669017
$ws.Range("D25").Value = -42686.75091417557
$ws.Range("E25").Value = 546465.3236527261
$ws.Range("B26").Value = 657658.1297490995
$ws.Range("D26").Value = -43483.74228148611
$ws.Range("E26").Value = 614174.3874676134
$ws.Range("B27").Value = 726164.1849312974
$ws.Range("D27").Value = -44066.51858443201
$ws.Range("E27").Value = 682097.6663468655

# --- Sheet: NZU ---
$ws = $wb.Worksheets.Item("NZU")
$ws.Range("B2").Value = -1059392.721239566
$ws.Range("D2").Value = 612422.837663063
$ws.Range("E2").Value = -446969.8835765031
$ws.Range("B3").Value = -1016152.202005298
$ws.Range("D3").Value = 571902.862059198
$ws.Range("E3").Value = -444249.3399461001
$ws.Range("B4").Value = -972911.6827710301
$ws.Range("D4").Value = 532136.861338568
$ws.Range("E4").Value = -440774.8214324621
$ws.Range("B5").Value = -929671.1635367621
$ws.Range("D5").Value = 493242.497808507
$ws.Range("E5").Value = -436428.6657282551
$ws.Range("B6").Value = -886430.6443024941
$ws.Range("D6").Value = 455336.2644986401
$ws.Range("E6").Value = -431094.379803854
$ws.Range("B7").Value = -843190.1250682261
$ws.Range("D7").Value = 418530.2743451419
$ws.Range("E7").Value = -424659.8507230842
$ws.Range("B8").Value = -799949.6058339581
$ws.Range("D8").Value = 382929.2781403388
$ws.Range("E8").Value = -417020.3276936193
$ws.Range("B9").Value = -756709.08659969
$ws.Range("D9").Value = 348628.051261787
$ws.Range("E9").Value = -408081.035337903
$ws.Range("B10").Value = -713468.5673654221
$ws.Range("D10").Value = 315709.2502861379
$ws.Range("E10").Value = -397759.3170792843
$ws.Range("B11").Value = -670228.048131154
$ws.Range("D11").Value = 284241.8000348282
$ws.Range("E11").Value = -385986.2480963258
$ws.Range("B12").Value = -626987.5288968862
$ws.Range("D12").Value = 254279.8331757253
$ws.Range("E12").Value = -372707.6957211609
$ws.Range("B13").Value = -583747.009662618
$ws.Range("D13").Value = 225862.1709402024
$ws.Range("E13").Value = -357884.8387224157
$ws.Range("B14").Value = -540506.4904283501
$ws.Range("D14").Value = 199012.3064569125
$ws.Range("E14").Value = -341494.1839714375
$ws.Range("B15").Value = -497265.9711940821
$ws.Range("D15").Value = 173738.8323405264
$ws.Range("E15").Value = -323527.1388535557
$ws.Range("B16").Value = -454025.4519598141
$ws.Range("D16").Value = 150036.2414234223
$ws.Range("E16").Value = -303989.2105363918
$ws.Range("B17").Value = -410784.9327255461
$ws.Range("D17").Value = 127886.0232417857
$ws.Range("E17").Value = -282898.9094837603
$ws.Range("B18").Value = -367544.4134912781
$ws.Range("D18").Value = 107257.9780972471
$ws.Range("E18").Value = -260286.435394031
$ws.Range("B19").Value = -324303.89425701
$ws.Range("D19").Value = 88111.67405544764
$ws.Range("E19").Value = -236192.2202015624
$ws.Range("B20").Value = -281063.375022742
$ws.Range("D20").Value = 70397.9789323185
$ws.Range("E20").Value = -210665.3960904235
$ws.Range("B21").Value = -237822.855788474
$ws.Range("D21").Value = 54060.60805210926
$ws.Range("E21").Value = -183762.2477363648
$ws.Range("B22").Value = -194582.336554206
$ws.Range("D22").Value = 39037.63837386193
$ws.Range("E22").Value = -155544.6981803441
$ws.Range("B23").Value = -151341.817319938
$ws.Range("D23").Value = 25262.94968460235
$ws.Range("E23").Value = -126078.8676353357
$ws.Range("B24").Value = -108101.29808567
$ws.Range("D24").Value = 12667.56334026144
$ws.Range("E24").Value = -95433.73474540855
$ws.Range("B25").Value = -64860.77885140201
$ws.Range("D25").Value = 1180.858065976813
$ws.Range("E25").Value = -63679.9207854252
$ws.Range("B26").Value = -21620.259617134
$ws.Range("D26").Value = -9268.349672258742
$ws.Range("E26").Value = -30888.60928939274
$ws.Range("B27").Value = 21620.259617134
$ws.Range("D27").Value = -18750.86639078331
$ws.Range("E27").Value = 2869.393226350694
$ws.Range("B28").Value = 64860.77885140201
$ws.Range("D28").Value = -27336.32347405757
$ws.Range("E28").Value = 37524.45537734444
$ws.Range("B29").Value = 108101.29808567
$ws.Range("D29").Value = -35092.51588380117
$ws.Range("E29").Value = 73008.78220186883
$ws.Range("B30").Value = 151341.817319938
$ws.Range("D30").Value = -42084.87005080884
$ws.Range("E30").Value = 109256.9472691292
$ws.Range("B31").Value = 194582.336554206
$ws.Range("D31").Value = -48376.03099362033
$ws.Range("E31").Value = 146206.3055605857
$ws.Range("B32").Value = 237822.855788474
$ws.Range("D32").Value = -54025.5572423669
$ws.Range("E32").Value = 183797.2985461071
$ws.Range("B33").Value = 281063.375022742
$ws.Range("D33").Value = -59089.71138755204
$ws.Range("E33").Value = 221973.66363519
$ws.Range("B34").Value = 324303.89425701
$ws.Range("D34").Value = -63621.33388403049
$ws.Range("E34").Value = 260682.5603729795
$ws.Range("B35").Value = 367544.4134912781
$ws.Range("D35").Value = -67669.78799192734
$ws.Range("E35").Value = 299874.6254993507
$ws.Range("B36").Value = 410784.9327255461
$ws.Range("D36").Value = -71280.96431412169
$ws.Range("E36").Value = 339503.9684114244
$ws.Range("B37").Value = 454025.4519598141
$ws.Range("D37").Value = -74497.33419534676
$ws.Range("E37").Value = 379528.1177644673
$ws.Range("B38").Value = 497265.9711940821
$ws.Range("D38").Value = -77358.04219793518
$ws.Range("E38").Value = 419907.9289961469
$ws.Range("B39").Value = 540506.4904283501
$ws.Range("D39").Value = -79899.028896271
$ws.Range("E39").Value = 460607.461532079
$ws.Range("B40").Value = 583747.009662618
$ws.Range("D40").Value = -82153.17628294621
$ws.Range("E40").Value = 501593.8333796719
$ws.Range("B41").Value = 626987.5288968862
$ws.Range("D41").Value = -84150.46911387723
$ws.Range("E41").Value = 542837.0597830089
$ws.Range("B42").Value = 670228.048131154
$ws.Range("D42").Value = -85918.16650783984
$ws.Range("E42").Value = 584309.8816233142
$ws.Range("B43").Value = 713468.5673654221
$ws.Range("D43").Value = -87480.97903767644
$ws.Range("E43").Value = 625987.5883277457
$ws.Range("B44").Value = 756709.08659969
$ws.Range("D44").Value = -88861.24739320669
$ws.Range("E44").Value = 667847.8392064833
$ws.Range("B45").Value = 799949.6058339581
$ws.Range("D45").Value = -90079.11945290214
$ws.Range("E45").Value = 709870.486381056
$ws.Range("B46").Value = 843190.1250682261
$ws.Range("D46").Value = -91152.72327067491
$ws.Range("E46").Value = 752037.4017975512
$ws.Range("B47").Value = 886430.6443024941
$ws.Range("D47").Value = -92098.33406710562
$ws.Range("E47").Value = 794332.3102353885
$ws.Range("B48").Value = 929671.1635367621
$ws.Range("D48").Value = -92930.53381502797
$ws.Range("E48").Value = 836740.6297217342
$ws.Range("B49").Value = 972911.6827710301
$ws.Range("D49").Value = -93662.36243317055
$ws.Range("E49").Value = 879249.3203378596
$ws.Range("B50").Value = 1016152.202005298
$ws.Range("D50").Value = -94305.45995515988
$ws.Range("E50").Value = 921846.7420501382
$ws.Range("B51").Value = 1059392.721239566
$ws.Range("D51").Value = -94870.19933169798
$ws.Range("E51").Value = 964522.5219078681
$ws.Range("B52").Value = 1102633.240473834
$ws.Range("D52").Value = -95365.80975841579
$ws.Range("E52").Value = 1007267.430715418
$ws.Range("B53").Value = 1145873.759708102
$ws.Range("D53").Value = -95800.49060782605
$ws.Range("E53").Value = 1050073.269100276
$ws.Range("B54").Value = 1189114.27894237
$ws.Range("D54").Value = -96181.51618774552
$ws.Range("E54").Value = 1092932.762754624
$ws.Range("B55").Value = 1232354.798176638
$ws.Range("D55").Value = -96515.33165671221
$ws.Range("E55").Value = 1135839.466519926
$ws.Range("B56").Value = 1275595.317410906
$ws.Range("D56").Value = -96807.64050509584
$ws.Range("E56").Value = 1178787.67690581
$ws.Range("B57").Value = 1318835.836645174
$ws.Range("D57").Value = -97063.48406381189
$ws.Range("E57").Value = 1221772.352581362
$ws.Range("B58").Value = 1362076.355879442
$ws.Range("D58").Value = -97287.31353540525
$ws.Range("E58").Value = 1264789.042344037
$ws.Range("B59").Value = 1405316.87511371
$ws.Range("D59").Value = -97483.05505868763
$ws.Range("E59").Value = 1307833.820055022
$ws.Range("B60").Value = 1448557.394347978
$ws.Range("D60").Value = -97654.1683215073
$ws.Range("E60").Value = 1350903.226026471
$ws.Range("B61").Value = 1491797.913582246
$ws.Range("D61").Value = -97803.69922952542
$ws.Range("E61").Value = 1393994.214352721
$ws.Range("B62").Value = 1535038.432816514
$ws.Range("D62").Value = -97934.32712450752
$ws.Range("E62").Value = 1437104.105692007

# --- Sheet: EUA ---
$ws = $wb.Worksheets.Item("EUA")
$ws.Range("B2").Value = -1354932.909375169
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = -1354932.909375169
$ws.Range("B3").Value = -1321543.729745422
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = -1321543.729745422
$ws.Range("B4").Value = -1288154.550115674
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -1288154.550115674
$ws.Range("B5").Value = -1254765.370485926
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = -1254765.370485926
$ws.Range("B6").Value = -1221376.190856178
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -1221376.190856178
$ws.Range("B7").Value = -1187987.01122643
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = -1187987.01122643
$ws.Range("B8").Value = -1154597.831596682
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = -1154597.831596682
$ws.Range("B9").Value = -1121208.651966934
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = -1121208.651966934
$ws.Range("B10").Value = -1087819.472337186
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -1087819.472337186
$ws.Range("B11").Value = -1054430.292707438
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = -1054430.292707438
$ws.Range("B12").Value = -1021041.11307769
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = -1021041.11307769
$ws.Range("B13").Value = -987651.9334479426
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = -987651.9334479426
$ws.Range("B14").Value = -954262.7538181947
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = -954262.7538181947
$ws.Range("B15").Value = -920873.5741884467
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = -920873.5741884467
$ws.Range("B16").Value = -887484.3945586988
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = -887484.3945586988
$ws.Range("B17").Value = -854095.214928951
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = -854095.214928951
$ws.Range("B18").Value = -820706.0352992031
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = -820706.0352992031
$ws.Range("B19").Value = -787316.8556694551
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = -787316.8556694551
$ws.Range("B20").Value = -753927.6760397073
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = -753927.6760397073
$ws.Range("B21").Value = -720538.4964099595
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = -720538.4964099595
$ws.Range("B22").Value = -687149.3167802115
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = -687149.3167802115
$ws.Range("B23").Value = -653760.1371504636
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = -653760.1371504636
$ws.Range("B24").Value = -620370.9575207158
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = -620370.9575207158
$ws.Range("B25").Value = -586981.7778909679
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = -586981.7778909679
$ws.Range("B26").Value = -553592.5982612199
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = -553592.5982612199
$ws.Range("B27").Value = -520203.4186314721
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = -520203.4186314721
$ws.Range("B28").Value = -486814.2390017241
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = -486814.2390017241
$ws.Range("B29").Value = -453425.0593719763
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = -453425.0593719763
$ws.Range("B30").Value = -420035.8797422284
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = -420035.8797422284
$ws.Range("B31").Value = -386646.7001124805
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = -386646.7001124805
$ws.Range("B32").Value = -353257.5204827326
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = -353257.5204827326
$ws.Range("B33").Value = -319868.3408529848
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = -319868.3408529848
$ws.Range("B34").Value = -286479.1612232369
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = -286479.1612232369
$ws.Range("B35").Value = -253089.9815934889
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = -253089.9815934889
$ws.Range("B36").Value = -219700.8019637411
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = -219700.8019637411
$ws.Range("B37").Value = -186311.6223339932
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = -186311.6223339932
$ws.Range("B38").Value = -152922.4427042453
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = -152922.4427042453
$ws.Range("B39").Value = -119533.2630744974
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = -119533.2630744974
$ws.Range("B40").Value = -86144.0834447495
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = -86144.0834447495
$ws.Range("B41").Value = -52754.90381500161
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = -52754.90381500161
$ws.Range("B42").Value = -19365.72418525372
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = -19365.72418525372
$ws.Range("B43").Value = 14023.45544449417
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 14023.45544449417
$ws.Range("B44").Value = 47412.63507424206
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 47412.63507424206
$ws.Range("B45").Value = 80801.81470398995
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 80801.81470398995
$ws.Range("B46").Value = 114190.9943337378
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 114190.9943337378
$ws.Range("B47").Value = 147580.1739634857
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 147580.1739634857
$ws.Range("B48").Value = 180969.3535932336
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 180969.3535932336
$ws.Range("B49").Value = 214358.5332229815
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 214358.5332229815
$ws.Range("B50").Value = 247747.7128527294
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 247747.7128527294
$ws.Range("B51").Value = 281136.8924824773
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 281136.8924824773
$ws.Range("B52").Value = 314526.0721122252
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 314526.0721122252
$ws.Range("B53").Value = 347915.2517419731
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 347915.2517419731
$ws.Range("B54").Value = 381304.431371721
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 381304.431371721
$ws.Range("B55").Value = 414693.6110014688
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 414693.6110014688
$ws.Range("B56").Value = 448082.7906312168
$ws.Range("D56").Value = 0
$ws.Range("E56").Value = 448082.7906312168
$ws.Range("B57").Value = 481471.9702609646
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 481471.9702609646
$ws.Range("B58").Value = 514861.1498907126
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 514861.1498907126
$ws.Range("B59").Value = 548250.3295204605
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 548250.3295204605
$ws.Range("B60").Value = 581639.5091502083
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 581639.5091502083
$ws.Range("B61").Value = 615028.6887799562
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 615028.6887799562
$ws.Range("B62").Value = 648417.8684097041
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 648417.8684097041
$ws.Range("B63").Value = 681807.0480394519
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 681807.0480394519
$ws.Range("B64").Value = 715196.2276691998
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 715196.2276691998
$ws.Range("B65").Value = 748585.4072989478
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 748585.4072989478
$ws.Range("B66").Value = 781974.5869286957
$ws.Range("D66").Value = 0
$ws.Range("E66").Value = 781974.5869286957
$ws.Range("B67").Value = 815363.7665584435
$ws.Range("D67").Value = 0
$ws.Range("E67").Value = 815363.7665584435
$ws.Range("B68").Value = 848752.9461881914
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 848752.9461881914
$ws.Range("B69").Value = 882142.1258179394
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 882142.1258179394
$ws.Range("B70").Value = 915531.3054476872
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 915531.3054476872
$ws.Range("B71").Value = 948920.485077435
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 948920.485077435
$ws.Range("B72").Value = 982309.664707183
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 982309.664707183

# --- Sheet: UKA ---
$ws = $wb.Worksheets.Item("UKA")
$ws.Range("B2").Value = -838957.4106661081
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = -838957.4106661081
$ws.Range("B3").Value = -778251.3751620481
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = -778251.3751620481
$ws.Range("B4").Value = -717545.3396579883
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -717545.3396579883
$ws.Range("B5").Value = -656839.3041539284
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = -656839.3041539284
$ws.Range("B6").Value = -596133.2686498684
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -596133.2686498684
$ws.Range("B7").Value = -535427.2331458085
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = -535427.2331458085
$ws.Range("B8").Value = -474721.1976417486
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = -474721.1976417486
$ws.Range("B9").Value = -414015.1621376887
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = -414015.1621376887
$ws.Range("B10").Value = -353309.1266336287
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -353309.1266336287
$ws.Range("B11").Value = -292603.0911295689
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = -292603.0911295689
$ws.Range("B12").Value = -231897.0556255089
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = -231897.0556255089
$ws.Range("B13").Value = -171191.020121449
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = -171191.020121449
$ws.Range("B14").Value = -110484.9846173891
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = -110484.9846173891
$ws.Range("B15").Value = -49778.94911332915
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = -49778.94911332915
$ws.Range("B16").Value = 10927.08639073077
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 10927.08639073077
$ws.Range("B17").Value = 71633.12189479069
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 71633.12189479069
$ws.Range("B18").Value = 132339.1573988506
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 132339.1573988506
$ws.Range("B19").Value = 193045.1929029105
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 193045.1929029105
$ws.Range("B20").Value = 253751.2284069704
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 253751.2284069704
$ws.Range("B21").Value = 314457.2639110303
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 314457.2639110303
$ws.Range("B22").Value = 375163.2994150903
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 375163.2994150903
$ws.Range("B23").Value = 435869.3349191502
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 435869.3349191502
$ws.Range("B24").Value = 496575.3704232101
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 496575.3704232101
$ws.Range("B25").Value = 557281.4059272701
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 557281.4059272701
$ws.Range("B26").Value = 617987.4414313299
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 617987.4414313299
$ws.Range("B27").Value = 678693.4769353899
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 678693.4769353899
$ws.Range("B28").Value = 739399.5124394498
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 739399.5124394498
$ws.Range("B29").Value = 800105.5479435096
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 800105.5479435096
$ws.Range("B30").Value = 860811.5834475696
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 860811.5834475696
$ws.Range("B31").Value = 921517.6189516296
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 921517.6189516296
$ws.Range("B32").Value = 982223.6544556894
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 982223.6544556894
$ws.Range("B33").Value = 1042929.689959749
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 1042929.689959749
$ws.Range("B34").Value = 1103635.725463809
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 1103635.725463809
$ws.Range("B35").Value = 1164341.760967869
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 1164341.760967869
$ws.Range("B36").Value = 1225047.796471929
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 1225047.796471929
$ws.Range("B37").Value = 1285753.831975989
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 1285753.831975989
$ws.Range("B38").Value = 1346459.867480049
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 1346459.867480049
$ws.Range("B39").Value = 1407165.902984109
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 1407165.902984109
$ws.Range("B40").Value = 1467871.938488169
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 1467871.938488169
$ws.Range("B41").Value = 1528577.973992229
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 1528577.973992229
$ws.Range("B42").Value = 1589284.009496289
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 1589284.009496289
$ws.Range("B43").Value = 1649990.045000348
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 1649990.045000348
$ws.Range("B44").Value = 1710696.080504408
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 1710696.080504408
$ws.Range("B45").Value = 1771402.116008468
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 1771402.116008468
$ws.Range("B46").Value = 1832108.151512528
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 1832108.151512528
$ws.Range("B47").Value = 1892814.187016588
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 1892814.187016588

# --- Sheet: CCA ---
$ws = $wb.Worksheets.Item("CCA")
$ws.Range("B2").Value = -637370.3534764132
$ws.Range("D2").Value = 310410.3746481487
$ws.Range("E2").Value = -326959.9788282645
$ws.Range("B3").Value = -569129.4162733712
$ws.Range("D3").Value = 260413.4498757628
$ws.Range("E3").Value = -308715.9663976084
$ws.Range("B4").Value = -500888.4790703291
$ws.Range("D4").Value = 213902.962343973
$ws.Range("E4").Value = -286985.5167263561
$ws.Range("B5").Value = -432647.541867287
$ws.Range("D5").Value = 171014.066250351
$ws.Range("E5").Value = -261633.475616936
$ws.Range("B6").Value = -364406.6046642449
$ws.Range("D6").Value = 131802.129831046
$ws.Range("E6").Value = -232604.4748331989
$ws.Range("B7").Value = -296165.6674612028
$ws.Range("D7").Value = 96247.50789685547
$ws.Range("E7").Value = -199918.1595643474
$ws.Range("B8").Value = -227924.7302581607
$ws.Range("D8").Value = 64264.04326782994
$ws.Range("E8").Value = -163660.6869903308
$ws.Range("B9").Value = -159683.7930551187
$ws.Range("D9").Value = 35709.93807499253
$ws.Range("E9").Value = -123973.8549801262
$ws.Range("B10").Value = -91442.85585207661
$ws.Range("D10").Value = 10399.7390515657
$ws.Range("E10").Value = -81043.11680051091
$ws.Range("B11").Value = -23201.91864903453
$ws.Range("D11").Value = -11883.60981592959
$ws.Range("E11").Value = -35085.52846496413
$ws.Range("B12").Value = 45039.01855400753
$ws.Range("D12").Value = -31377.43745777605
$ws.Range("E12").Value = 13661.58109623149
$ws.Range("B13").Value = 113279.9557570496
$ws.Range("D13").Value = -48329.26994060729
$ws.Range("E13").Value = 64950.6858164423
$ws.Range("B14").Value = 181520.8929600917
$ws.Range("D14").Value = -62988.30778214676
$ws.Range("E14").Value = 118532.5851779449
$ws.Range("B15").Value = 249761.8301631338
$ws.Range("D15").Value = -75598.56561138935
$ws.Range("E15").Value = 174163.2645517444
$ws.Range("B16").Value = 318002.7673661758
$ws.Range("D16").Value = -86393.6411115244
$ws.Range("E16").Value = 231609.1262546514
$ws.Range("B17").Value = 386243.7045692179
$ws.Range("D17").Value = -95592.9833165
$ws.Range("E17").Value = 290650.7212527179
$ws.Range("B18").Value = 454484.64177226
$ws.Range("D18").Value = -103399.4727408348
$ws.Range("E18").Value = 351085.1690314251
$ws.Range("B19").Value = 522725.578975302
$ws.Range("D19").Value = -109998.0995202294
$ws.Range("E19").Value = 412727.4794550726
$ws.Range("B20").Value = 590966.516178344
$ws.Range("D20").Value = -115555.5225301925
$ws.Range("E20").Value = 475410.9936481515
$ws.Range("B21").Value = 659207.4533813862
$ws.Range("D21").Value = -120220.3048952044
$ws.Range("E21").Value = 538987.1484861819
$ws.Range("B22").Value = 727448.3905844283
$ws.Range("D22").Value = -124123.6432125356
$ws.Range("E22").Value = 603324.7473718928
$ws.Range("B23").Value = 795689.3277874703
$ws.Range("D23").Value = -127380.4343953814
$ws.Range("E23").Value = 668308.8933920889
$ws.Range("B24").Value = 863930.2649905123
$ws.Range("D24").Value = -130090.5518170293
$ws.Range("E24").Value = 733839.713173483
$ws.Range("B25").Value = 932171.2021935545
$ws.Range("D25").Value = -132340.2290833753
$ws.Range("E25").Value = 799830.9731101792
$ws.Range("B26").Value = 1000412.139396597
$ws.Range("D26").Value = -134203.4738666353
$ws.Range("E26").Value = 866208.6655299612
$ws.Range("B27").Value = 1068653.076599639
$ws.Range("D27").Value = -135743.4550836609
$ws.Range("E27").Value = 932909.6215159778
$ws.Range("B28").Value = 1136894.013802681
$ws.Range("D28").Value = -137013.8240749925
$ws.Range("E28").Value = 999880.1897276882
$ws.Range("B29").Value = 1205134.951005723
$ws.Range("D29").Value = -138059.9444377709
$ws.Range("E29").Value = 1067075.006567952
$ws.Range("B30").Value = 1273375.888208765
$ws.Range("D30").Value = -138920.0160841833
$ws.Range("E30").Value = 1134455.872124582
$ws.Range("B31").Value = 1341616.825411807
$ws.Range("D31").Value = -139626.0873341532
$ws.Range("E31").Value = 1201990.738077654
$ws.Range("B32").Value = 1409857.762614849
$ws.Range("D32").Value = -140204.9548366376
$ws.Range("E32").Value = 1269652.807778212
$ws.Range("B33").Value = 1478098.699817891
$ws.Range("D33").Value = -140678.9552683685
$ws.Range("E33").Value = 1337419.744549523
$ws.Range("B34").Value = 1546339.637020933
$ws.Range("D34").Value = -141066.6554668855
$ws.Range("E34").Value = 1405272.981554048
$ws.Range("B35").Value = 1614580.574223975
$ws.Range("D35").Value = -141383.4492524743
$ws.Range("E35").Value = 1473197.124971501
$ws.Range("B36").Value = 1682821.511427017
$ws.Range("D36").Value = -141642.0699645763
$ws.Range("E36").Value = 1541179.441462441
